$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Set up header row and data row (write data row values first so the
# shared-string table keeps "Admin"/"admin123" at indices 0/1, matching
# the original "Akshara"/"Swara" slots, with "Username"/"Password"
# appended afterwards at indices 2/3)
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Update selection to A7:XFD7 (select entire row 7)
$ws.Rows(7).Select()
